$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("40:41").Delete()
$ws.Range("A2").Value = 'Multi-Utilities(18)'
$ws.Range("B2").Value = 0.6011530979219659
$ws.Range("A3").Value = 'Road & Rail(22)'
$ws.Range("B3").Value = 0.555478110834445
$ws.Range("A4").Value = 'Air Freight & Logistics(11)'
$ws.Range("B4").Value = 0.5120315985335757
$ws.Range("A5").Value = 'Electric Utilities(28)'
$ws.Range("B5").Value = 0.5111832501130005
$ws.Range("A6").Value = 'Banks(246)'
$ws.Range("B6").Value = 0.4914866227540425
$ws.Range("A7").Value = 'Building Products(23)'
$ws.Range("B7").Value = 0.4817696952719114
$ws.Range("A8").Value = 'Energy Equipment & Services(32)'
$ws.Range("B8").Value = 0.4718498402254432
$ws.Range("A9").Value = 'Metals & Mining(89)'
$ws.Range("B9").Value = 0.4316523242605258
$ws.Range("A10").Value = 'Machinery(85)'
$ws.Range("B10").Value = 0.4152030506263252
$ws.Range("A11").Value = 'Marine(15)'
$ws.Range("B11").Value = 0.3878966859677347
$ws.Range("A12").Value = 'Gas Utilities(12)'
$ws.Range("B12").Value = 0.3704453840804315
$ws.Range("A13").Value = 'Trading Companies & Distributors(25)'
$ws.Range("B13").Value = 0.3654504162682767
$ws.Range("A14").Value = 'Thrifts & Mortgage Finance(47)'
$ws.Range("B14").Value = 0.3547728109197359
$ws.Range("A15").Value = 'Water Utilities(12)'
$ws.Range("B15").Value = 0.3520959880163748
$ws.Range("A16").Value = 'Insurance(75)'
$ws.Range("B16").Value = 0.3180148923140174
$ws.Range("A17").Value = 'Auto Components(21)'
$ws.Range("B17").Value = 0.3055511741165613
$ws.Range("A18").Value = 'Construction & Engineering(20)'
$ws.Range("B18").Value = 0.3023106614107945
$ws.Range("A19").Value = 'Chemicals(51)'
$ws.Range("B19").Value = 0.2921011422361156
$ws.Range("A20").Value = 'Life Sciences Tools & Services(19)'
$ws.Range("B20").Value = 0.2631741957874983
$ws.Range("A21").Value = 'Specialty Retail(58)'
$ws.Range("B21").Value = 0.2409116727380363
$ws.Range("A22").Value = 'Capital Markets(75)'
$ws.Range("B22").Value = 0.2380249774285201
$ws.Range("A23").Value = 'Semiconductors & Semiconductor Equipment(68)'
$ws.Range("B23").Value = 0.2154474925683811
$ws.Range("A24").Value = 'Electrical Equipment(28)'
$ws.Range("B24").Value = 0.1998757633264295
$ws.Range("A25").Value = 'Commercial Services & Supplies(52)'
$ws.Range("B25").Value = 0.1914991355344332
$ws.Range("A26").Value = 'Professional Services(35)'
$ws.Range("B26").Value = 0.1896722757203425
$ws.Range("A27").Value = 'Aerospace & Defense(37)'
$ws.Range("B27").Value = 0.1828968893596951
$ws.Range("A28").Value = 'Hotels, Restaurants & Leisure(50)'
$ws.Range("B28").Value = 0.1787714773155454
$ws.Range("A29").Value = 'Oil, Gas & Consumable Fuels(122)'
$ws.Range("B29").Value = 0.1783105191961414
$ws.Range("A30").Value = 'Pharmaceuticals(48)'
$ws.Range("B30").Value = 0.1745937056713801
$ws.Range("A31").Value = 'Health Care Providers & Services(46)'
$ws.Range("B31").Value = 0.1591865300992264
$ws.Range("A32").Value = 'Communications Equipment(45)'
$ws.Range("B32").Value = 0.1519806083713743
$ws.Range("A33").Value = 'Media(42)'
$ws.Range("B33").Value = 0.1391667115132491
$ws.Range("A34").Value = 'Household Durables(39)'
$ws.Range("B34").Value = 0.1349024315186565
$ws.Range("A35").Value = 'Health Care Equipment & Supplies(83)'
$ws.Range("B35").Value = 0.1174498816020266
$ws.Range("A36").Value = 'Biotechnology(126)'
$ws.Range("B36").Value = 0.1129795663768614
$ws.Range("A37").Value = 'Food Products(44)'
$ws.Range("B37").Value = 0.1015430993239691
$ws.Range("A38").Value = 'IT Services(52)'
$ws.Range("B38").Value = 0.09920933229872297
$ws.Range("A39").Value = 'Software(66)'
$ws.Range("B39").Value = 0.08997060194791422
